# Update the "dissertation progress" tracker:
#  - Word count (B2) moves from 15560 to 16242 (dependent formulas recalc automatically)
#  - The logged snapshot value in B11 moves from 14697 to 16192, and its cell
#    border is cleared (it no longer carries the boxed/bordered style)
#  - The last active selection moves from D13 to E17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Total words written so far
$ws.Range("B2").Value = 16242

# Latest count snapshot, and drop the border that used to box this cell
$ws.Range("B11").Value = 16192
$ws.Range("B11").Borders.LineStyle = -4142   # xlLineStyleNone

# Leave the selection where the author left off
[void]$ws.Range("E17").Select()
